$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E4").Value = 3
$ws.Range("E18").Value = 4
$ws.Range("E28").Value = 4
$ws.Range("E29").Value = 2
$ws.Range("E32").Value = 3
$ws.Range("E38").Value = 2
$ws.Range("E39").Value = 4
$ws.Range("E46").Value = 6
$ws.Range("E60").Value = 3
$ws.Range("E65").Value = 5
$ws.Range("E68").Value = 2
$ws.Range("E83").Value = 3
$ws.Range("E113").Value = 3
$ws.Range("E147").Value = 3
$ws.Range("E164").Value = 3
$ws.Range("E176").Value = 3
$ws.Range("E177").Value = 4
$ws.Range("E180").Value = 3
$ws.Range("E190").Value = 2
$ws.Range("E197").Value = 2
$ws.Range("E231").Value = 2
$ws.Range("E297").Value = 2
$ws.Range("E305").Value = 2
$ws.Range("E311").Value = 3
$ws.Range("E314").Value = 2
$ws.Range("E328").Value = 6
$ws.Range("E329").Value = 4
$ws.Range("E343").Value = 4
$ws.Range("E365").Value = 4
$ws.Range("E398").Value = 2
$ws.Range("E408").Value = 2
$ws.Range("E440").Value = 3
$ws.Range("E442").Value = 3
$ws.Range("E443").Value = 4
$ws.Range("E450").Value = 3
$ws.Range("E467").Value = 5
$ws.Range("E477").Value = 2
$ws.Range("E483").Value = 2
$ws.Range("E530").Value = 5
$ws.Range("E542").Value = 4
$ws.Range("E560").Value = 2
$ws.Range("E562").Value = 2
$ws.Range("E564").Value = 4
$ws.Range("E571").Value = 4
$ws.Range("E577").Value = 2
$ws.Range("E585").Value = 3
$ws.Range("E613").Value = 4
$ws.Range("E617").Value = 2
$ws.Range("E639").Value = 3
$ws.Range("E645").Value = 2
$ws.Range("E649").Value = 2
$ws.Range("E659").Value = 4
$ws.Range("E673").Value = 2
$ws.Range("E677").Value = 3
$ws.Range("E678").Value = 4
$ws.Range("E690").Value = 4
$ws.Range("E717").Value = 2
$ws.Range("E729").Value = 5
$ws.Range("E732").Value = 2
$ws.Range("E756").Value = 2
$ws.Range("E772").Value = 4
$ws.Range("E775").Value = 2
$ws.Range("E790").Value = 3
$ws.Range("E795").Value = 3
$ws.Range("E812").Value = 3
$ws.Range("E829").Value = 2
$ws.Range("E838").Value = 2
$ws.Range("E840").Value = 3
$ws.Range("E841").Value = 5
$ws.Range("E848").Value = 2
$ws.Range("E871").Value = 4
$ws.Range("E877").Value = 2
$ws.Range("E885").Value = 3
$ws.Range("E901").Value = 3
$ws.Range("E902").Value = 4
$ws.Range("E907").Value = 3
$ws.Range("E908").Value = 2
$ws.Range("E912").Value = 2
$ws.Range("E918").Value = 2
$ws.Range("E935").Value = 6
$ws.Range("E947").Value = 3
$ws.Range("E953").Value = 2
$ws.Range("E975").Value = 2
$ws.Range("E986").Value = 2
$ws.Range("E1014").Value = 2
$ws.Range("E1070").Value = 4
$ws.Range("E1083").Value = 2
$ws.Range("E1084").Value = 2
$ws.Range("E1093").Value = 2
$ws.Range("E1095").Value = 2
$ws.Range("E1100").Value = 3
$ws.Range("E1118").Value = 3
$ws.Range("E1119").Value = 4
$ws.Range("E1123").Value = 2
$ws.Range("E1144").Value = 2
$ws.Range("E1159").Value = 2
$ws.Range("E1187").Value = 2
$ws.Range("E1188").Value = 4
$ws.Range("E1194").Value = 2
$ws.Range("E1203").Value = 2
$ws.Range("E1223").Value = 2
$ws.Range("E1224").Value = 4
$ws.Range("E1246").Value = 3
$ws.Range("E1263").Value = 3
$ws.Range("E1264").Value = 5
$ws.Range("E1275").Value = 3
$ws.Range("E1279").Value = 2
$ws.Range("E1283").Value = 2
$ws.Range("E1290").Value = 4
$ws.Range("E1299").Value = 2
$ws.Range("E1301").Value = 3
$ws.Range("E1306").Value = 2
